$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New version of the contact-tracing edge list data.
# Same graph topology as before, but node identifiers were regenerated.
$data = @(
    @("from",   "to"),
    @("51883d", "185911"),
    @("b4d8aa", "e4b0a2"),
    @("39e9dc", "b4d8aa"),
    @("39e9dc", "601d2e"),
    @("51883d", "9aa197"),
    @("39e9dc", "51883d"),
    @("39e9dc", "e399b1"),
    @("b4d8aa", "af0ac0"),
    @("39e9dc", "947e40"),
    @("39e9dc", "664549"),
    @("39e9dc", "605322")
)

# Values that look like plain numbers need the cell pre-formatted as
# Text, otherwise Excel would silently coerce them into numeric cells.
$numericLooking = @("185911", "947e40", "664549", "605322")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)

    if ($numericLooking -contains $data[$i][0]) {
        $cellA.NumberFormat = "@"
    }
    if ($numericLooking -contains $data[$i][1]) {
        $cellB.NumberFormat = "@"
    }

    $cellA.Value = $data[$i][0]
    $cellB.Value = $data[$i][1]
}
